$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point drift on the existing row 3 timestamp
# (re-saved by the scheduled task with slightly different precision).
$ws.Range("A3").Value = 45863.41686092593

# Append the new row 4 produced by the scheduled data-collection task.
$ws.Range("A4").Value = 45863.45853578657
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat

$ws.Range("B4").Value = 2025
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 17.67
$ws.Range("E4").Value = 78.55
$ws.Range("F4").Value = 571.29
$ws.Range("G4").Value = 12.52
$ws.Range("H4").Value = "ESE"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "11:00:17"
